$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Header text: issue/volume number and reporting week dates updated
# ------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  44"
$ws.Range("C9").Value = "Report Covering the Week  10/28/2024  Through  11/3/2024"

# ------------------------------------------------------------------
# Straight numeric value updates (stat tables, no type change)
# ------------------------------------------------------------------
$ws.Range("N14").Value = -93.75
$ws.Range("F15").Value = 5
$ws.Range("H15").Value = 400
$ws.Range("I15").Value = 19
$ws.Range("K15").Value = 46.153846153846
$ws.Range("L15").Value = 72.727272727272
$ws.Range("M15").Value = -36.666666666666
$ws.Range("N15").Value = -62
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 130
$ws.Range("K16").Value = -17.692307692307
$ws.Range("L16").Value = 15.053763440860
$ws.Range("M16").Value = -56.680161943319
$ws.Range("N16").Value = -87.396937573616
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 11
$ws.Range("E17").Value = -45.454545454545
$ws.Range("F17").Value = 39
$ws.Range("G17").Value = 40
$ws.Range("H17").Value = -2.5
$ws.Range("I17").Value = 412
$ws.Range("J17").Value = 399
$ws.Range("K17").Value = 3.258145363408
$ws.Range("L17").Value = 22.985074626865
$ws.Range("M17").Value = 59.073359073359
$ws.Range("N17").Value = -38.507462686567
$ws.Range("C18").Value = 2
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -64.285714285714
$ws.Range("I18").Value = 89
$ws.Range("J18").Value = 110
$ws.Range("K18").Value = -19.090909090909
$ws.Range("L18").Value = -4.301075268817
$ws.Range("M18").Value = -66.030534351145
$ws.Range("N18").Value = -93.933197000681
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 27
$ws.Range("G19").Value = 24
$ws.Range("H19").Value = 12.5
$ws.Range("I19").Value = 322
$ws.Range("J19").Value = 356
$ws.Range("K19").Value = -9.550561797752
$ws.Range("L19").Value = 11.418685121107
$ws.Range("M19").Value = -6.936416184971
$ws.Range("N19").Value = -34.285714285714
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -66.666666666666
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = -42.857142857142
$ws.Range("I20").Value = 111
$ws.Range("J20").Value = 134
$ws.Range("K20").Value = -17.164179104477
$ws.Range("L20").Value = 11
$ws.Range("M20").Value = -34.705882352941
$ws.Range("N20").Value = -90.176991150442
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = -20
$ws.Range("F21").Value = 94
$ws.Range("G21").Value = 103
$ws.Range("H21").Value = -8.737864077669
$ws.Range("I21").Value = 1061
$ws.Range("J21").Value = 1154
$ws.Range("K21").Value = -8.058925476603
$ws.Range("L21").Value = 14.455231930960
$ws.Range("M21").Value = -19.924528301886
$ws.Range("N21").Value = -77.290239726027
$ws.Range("F23").Value = 6
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 70
$ws.Range("K23").Value = -2.777777777777
$ws.Range("L23").Value = 42.857142857142
$ws.Range("M23").Value = 79.487179487179
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 34
$ws.Range("E24").Value = -50
$ws.Range("F24").Value = 93
$ws.Range("G24").Value = 108
$ws.Range("H24").Value = -13.888888888888
$ws.Range("I24").Value = 1029
$ws.Range("J24").Value = 1066
$ws.Range("K24").Value = -3.470919324577
$ws.Range("L24").Value = 2.9
$ws.Range("M24").Value = -20.417633410672
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 35
$ws.Range("G25").Value = 44
$ws.Range("H25").Value = -20.454545454545
$ws.Range("I25").Value = 446
$ws.Range("J25").Value = 404
$ws.Range("K25").Value = 10.396039603960
$ws.Range("L25").Value = 38.080495356037
$ws.Range("C26").Value = 13
$ws.Range("D26").Value = 19
$ws.Range("E26").Value = -31.578947368421
$ws.Range("F26").Value = 72
$ws.Range("G26").Value = 65
$ws.Range("H26").Value = 10.769230769230
$ws.Range("I26").Value = 673
$ws.Range("J26").Value = 621
$ws.Range("K26").Value = 8.373590982286
$ws.Range("L26").Value = 25.794392523364
$ws.Range("M26").Value = -38.200183654729
$ws.Range("F27").Value = 5
$ws.Range("H27").Value = 400
$ws.Range("I27").Value = 33
$ws.Range("K27").Value = 73.684210526315
$ws.Range("L27").Value = 13.793103448275
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 9
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = 28.571428571428
$ws.Range("I28").Value = 76
$ws.Range("J28").Value = 78
$ws.Range("K28").Value = -2.564102564102
$ws.Range("L28").Value = 16.923076923076
$ws.Range("E29").Value = -100
$ws.Range("J29").Value = 22
$ws.Range("K29").Value = -54.545454545454
$ws.Range("N29").Value = -89.473684210526
$ws.Range("E30").Value = -100
$ws.Range("J30").Value = 20
$ws.Range("K30").Value = -50
$ws.Range("N30").Value = -87.951807228915
$ws.Range("J31").Value = 4
$ws.Range("K31").Value = 0

# ------------------------------------------------------------------
# Cells switching from a numeric value to a text placeholder ("0" or
# "***.*") because the prior-year comparison value became unavailable
# ------------------------------------------------------------------
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "0"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "***.*"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "0"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "***.*"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "***.*"

# ------------------------------------------------------------------
# Cells switching from a text placeholder to an actual numeric value
# ------------------------------------------------------------------
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("C15").Value = 4
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value = 4
$ws.Range("D31").NumberFormat = "#,##0"
$ws.Range("D31").Value = 1
$ws.Range("E31").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E31").Value = -100
$ws.Range("G31").NumberFormat = "#,##0"
$ws.Range("G31").Value = 1
$ws.Range("H31").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H31").Value = -100
